# Auto-generated edit script applying the Raiden_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 1713.1333
$ws.Range("J32").Value = 1877.6
$ws.Range("L32").Value = 1877.6
$ws.Range("N32").Value = -2529.6
$ws.Range("H33").Value = 6518.625
$ws.Range("I33").Value = 176.92308
$ws.Range("K33").Value = 176.92308
$ws.Range("M33").Value = 52.07692
$ws.Range("H51").Value = 5782.7334
$ws.Range("J51").Value = 6007.846
$ws.Range("L51").Value = 6007.846
$ws.Range("N51").Value = -6975.846
$ws.Range("H70").Value = 124676.16
$ws.Range("J70").Value = 8149.5
$ws.Range("L70").Value = 24448.5
$ws.Range("N70").Value = -24988.5
$ws.Range("H73").Value = 124676.16
$ws.Range("J73").Value = 8149.5
$ws.Range("L73").Value = 24448.5
$ws.Range("N73").Value = -26320.5
$ws.Range("H76").Value = 7534.45
$ws.Range("I76").Value = 7617.3076
$ws.Range("J76").Value = 7380.5713
$ws.Range("K76").Value = 7617.3076
$ws.Range("L76").Value = 7380.5713
$ws.Range("M76").Value = -7302.3076
$ws.Range("N76").Value = -8010.5713
$ws.Range("H79").Value = 7534.45
$ws.Range("I79").Value = 7617.3076
$ws.Range("J79").Value = 7380.5713
$ws.Range("K79").Value = 7617.3076
$ws.Range("L79").Value = 7380.5713
$ws.Range("M79").Value = -6525.3076
$ws.Range("N79").Value = -9564.5713
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 10000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -21232
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 673.3333
$ws.Range("I100").Value = 747.5
$ws.Range("J100").Value = 525
$ws.Range("K100").Value = 747.5
$ws.Range("L100").Value = 525
$ws.Range("M100").Value = -206.5
$ws.Range("N100").Value = -1607
$ws.Range("H103").Value = 1135.1904
$ws.Range("I103").Value = 684.46155
$ws.Range("J103").Value = 1867.625
$ws.Range("K103").Value = 2053.38465
$ws.Range("L103").Value = 5602.875
$ws.Range("M103").Value = -1467.38465
$ws.Range("N103").Value = -6774.875
$ws.Range("H111").Value = 2503.625
$ws.Range("I111").Value = 2503.625
$ws.Range("K111").Value = 7510.875
$ws.Range("M111").Value = -4443.875
$ws.Range("H112").Value = 1724.9333
$ws.Range("I112").Value = 1298.2222
$ws.Range("J112").Value = 2365
$ws.Range("K112").Value = 3894.6666
$ws.Range("L112").Value = 7095
$ws.Range("M112").Value = -2786.6666
$ws.Range("N112").Value = -9311
$ws.Range("H116").Value = 3352.7273
$ws.Range("I116").Value = 2931.111
$ws.Range("J116").Value = 5250
$ws.Range("K116").Value = 2931.111
$ws.Range("L116").Value = 5250
$ws.Range("M116").Value = 510.8890000000001
$ws.Range("N116").Value = -12134
$ws.Range("H125").Value = 9572.666999999999
$ws.Range("J125").Value = 4536
$ws.Range("L125").Value = 40824
$ws.Range("N125").Value = -45744
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H133").Value = 67249.75
$ws.Range("J133").Value = 67249.75
$ws.Range("L133").Value = 67249.75
$ws.Range("N133").Value = -77369.75
$ws.Range("H138").Value = 3148.889
$ws.Range("I138").Value = 1621.1111
$ws.Range("J138").Value = 3658.1482
$ws.Range("K138").Value = 4863.3333
$ws.Range("L138").Value = 10974.4446
$ws.Range("M138").Value = 276.6666999999998
$ws.Range("N138").Value = -21254.4446
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 12546
$ws.Range("I32").Value = 4516.6665
$ws.Range("K32").Value = 4516.6665
$ws.Range("M32").Value = -4229.6665
$ws.Range("H34").Value = 7716392.5
$ws.Range("I34").Value = 20015200
$ws.Range("K34").Value = 20015200
$ws.Range("M34").Value = -20014929
$ws.Range("H45").Value = 2104.111
$ws.Range("I45").Value = 3140.2
$ws.Range("K45").Value = 3140.2
$ws.Range("M45").Value = -2763.2
$ws.Range("H61").Value = 3136.5789
$ws.Range("J61").Value = 4809
$ws.Range("L61").Value = 4809
$ws.Range("N61").Value = -5233
$ws.Range("H62").Value = 30249
$ws.Range("J62").Value = 30249
$ws.Range("L62").Value = 30249
$ws.Range("N62").Value = -31497
$ws.Range("H65").Value = 30249
$ws.Range("J65").Value = 30249
$ws.Range("L65").Value = 90747
$ws.Range("N65").Value = -96987
$ws.Range("H74").Value = 976.3333
$ws.Range("I74").Value = 719.875
$ws.Range("J74").Value = 1489.25
$ws.Range("K74").Value = 719.875
$ws.Range("L74").Value = 1489.25
$ws.Range("M74").Value = 154.125
$ws.Range("N74").Value = -3237.25
$ws.Range("H77").Value = 976.3333
$ws.Range("I77").Value = 719.875
$ws.Range("J77").Value = 1489.25
$ws.Range("K77").Value = 3599.375
$ws.Range("L77").Value = 7446.25
$ws.Range("M77").Value = 768.625
$ws.Range("N77").Value = -16182.25
$ws.Range("H109").Value = 67598.60000000001
$ws.Range("J109").Value = 67598.60000000001
$ws.Range("L109").Value = 67598.60000000001
$ws.Range("N109").Value = -70372.60000000001
$ws.Range("H110").Value = 3271.125
$ws.Range("I110").Value = 2961.6667
$ws.Range("J110").Value = 4199.5
$ws.Range("K110").Value = 2961.6667
$ws.Range("L110").Value = 4199.5
$ws.Range("M110").Value = -916.6667000000002
$ws.Range("N110").Value = -8289.5
$ws.Range("H132").Value = 3362.9
$ws.Range("I132").Value = 3177.6072
$ws.Range("J132").Value = 5957
$ws.Range("K132").Value = 9532.821599999999
$ws.Range("L132").Value = 17871
$ws.Range("M132").Value = -7002.821599999999
$ws.Range("N132").Value = -22931
$ws.Range("H136").Value = 3136.5789
$ws.Range("J136").Value = 4809
$ws.Range("L136").Value = 14427
$ws.Range("N136").Value = -19527
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 1344.5714
$ws.Range("I20").Value = 1282.6
$ws.Range("K20").Value = 1282.6
$ws.Range("M20").Value = -1035.6
$ws.Range("H80").Value = 1605.2307
$ws.Range("I80").Value = 124.57143
$ws.Range("J80").Value = 3332.6667
$ws.Range("K80").Value = 124.57143
$ws.Range("L80").Value = 3332.6667
$ws.Range("M80").Value = 873.42857
$ws.Range("N80").Value = -5328.6667
$ws.Range("H83").Value = 1605.2307
$ws.Range("I83").Value = 124.57143
$ws.Range("J83").Value = 3332.6667
$ws.Range("K83").Value = 622.85715
$ws.Range("L83").Value = 16663.3335
$ws.Range("M83").Value = 4369.14285
$ws.Range("N83").Value = -26647.3335
$ws.Range("H94").Value = 1512.4642
$ws.Range("I94").Value = 1541.8889
$ws.Range("J94").Value = 1459.5
$ws.Range("K94").Value = 1541.8889
$ws.Range("L94").Value = 1459.5
$ws.Range("M94").Value = -1090.8889
$ws.Range("N94").Value = -2361.5
$ws.Range("H105").Value = 3545.353
$ws.Range("I105").Value = 2052.25
$ws.Range("J105").Value = 4004.7693
$ws.Range("K105").Value = 2052.25
$ws.Range("L105").Value = 4004.7693
$ws.Range("M105").Value = -305.25
$ws.Range("N105").Value = -7498.7693
$ws.Range("H134").Value = 4754.2856
$ws.Range("J134").Value = 4289.5
$ws.Range("L134").Value = 12868.5
$ws.Range("N134").Value = -17938.5
$ws = $wb.Worksheets.Item(4)
$ws.Range("H4").Value = 16169000
$ws.Range("J4").Value = 16169000
$ws.Range("L4").Value = 16169000
$ws.Range("N4").Value = -16169224
$ws.Range("H11").Value = 10000
$ws.Range("J11").Value = 10000
$ws.Range("L11").Value = 10000
$ws.Range("N11").Value = -10280
$ws.Range("H16").Value = 199.5
$ws.Range("I16").Value = 199.5
$ws.Range("K16").Value = 199.5
$ws.Range("M16").Value = 87.5
$ws.Range("H22").Value = 694.7273
$ws.Range("I22").Value = 632.3333
$ws.Range("J22").Value = 769.6
$ws.Range("K22").Value = 632.3333
$ws.Range("L22").Value = 769.6
$ws.Range("M22").Value = -282.3333
$ws.Range("N22").Value = -1469.6
$ws.Range("H31").Value = 7435.84
$ws.Range("I31").Value = 3398.9333
$ws.Range("J31").Value = 13491.2
$ws.Range("K31").Value = 3398.9333
$ws.Range("L31").Value = 13491.2
$ws.Range("M31").Value = -3103.9333
$ws.Range("N31").Value = -14081.2
$ws.Range("H34").Value = 7435.84
$ws.Range("I34").Value = 3398.9333
$ws.Range("J34").Value = 13491.2
$ws.Range("K34").Value = 3398.9333
$ws.Range("L34").Value = 13491.2
$ws.Range("M34").Value = -3196.9333
$ws.Range("N34").Value = -13895.2
$ws.Range("H35").Value = 202820.2
$ws.Range("I35").Value = 337666.66
$ws.Range("J35").Value = 550.5
$ws.Range("K35").Value = 337666.66
$ws.Range("L35").Value = 550.5
$ws.Range("M35").Value = -337372.66
$ws.Range("N35").Value = -1138.5
$ws.Range("H58").Value = 3065
$ws.Range("I58").Value = 3065
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3065
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2862
$ws.Range("N58").ClearContents()
$ws.Range("H68").Value = 73250
$ws.Range("J68").Value = 73250
$ws.Range("L68").Value = 73250
$ws.Range("N68").Value = -74748
$ws.Range("H71").Value = 73250
$ws.Range("J71").Value = 73250
$ws.Range("L71").Value = 219750
$ws.Range("N71").Value = -227238
$ws.Range("H107").Value = 648.4
$ws.Range("I107").Value = 414.33334
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 414.33334
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 1505.66666
$ws.Range("N107").Value = -4839.5
$ws.Range("H113").Value = 199.5
$ws.Range("I113").Value = 199.5
$ws.Range("K113").Value = 199.5
$ws.Range("M113").Value = 1970.5
$ws.Range("H122").Value = 2254.1428
$ws.Range("I122").Value = 2380
$ws.Range("K122").Value = 7140
$ws.Range("M122").Value = -4690
$ws.Range("H132").Value = 2989.25
$ws.Range("I132").Value = 2915.8333
$ws.Range("J132").Value = 3062.6667
$ws.Range("K132").Value = 8747.499899999999
$ws.Range("L132").Value = 9188.000100000001
$ws.Range("M132").Value = -6217.499899999999
$ws.Range("N132").Value = -14248.0001
$ws.Range("H135").Value = 90389
$ws.Range("J135").Value = 90389
$ws.Range("L135").Value = 90389
$ws.Range("N135").Value = -100529
$ws.Range("H136").Value = 3065
$ws.Range("I136").Value = 3065
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9195
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -6645
$ws.Range("N136").ClearContents()
$ws.Range("H141").Value = 291540.22
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 326732.75
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 326732.75
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -337092.75
$ws = $wb.Worksheets.Item(5)
$ws.Range("H11").Value = 469.13635
$ws.Range("I11").Value = 251.38889
$ws.Range("K11").Value = 754.1666700000001
$ws.Range("M11").Value = -614.1666700000001
$ws.Range("H80").Value = 2999
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10869
$ws.Range("H83").Value = 2999
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36351
$ws.Range("H94").Value = 2950
$ws.Range("I94").Value = 2950
$ws.Range("K94").Value = 8850
$ws.Range("M94").Value = -8174
$ws.Range("H98").Value = 479.66666
$ws.Range("I98").Value = 472.5
$ws.Range("K98").Value = 1417.5
$ws.Range("M98").Value = 80.5
$ws.Range("H109").Value = 3964.9375
$ws.Range("I109").Value = 3066.4
$ws.Range("K109").Value = 9199.200000000001
$ws.Range("M109").Value = -8159.200000000001
$ws.Range("H122").Value = 1673.7273
$ws.Range("I122").Value = 1398.8
$ws.Range("J122").Value = 1902.8334
$ws.Range("K122").Value = 12589.2
$ws.Range("L122").Value = 17125.5006
$ws.Range("M122").Value = -10139.2
$ws.Range("N122").Value = -22025.5006
$ws.Range("H131").Value = 54461.76
$ws.Range("I131").Value = 111930.7
$ws.Range("K131").Value = 335792.1
$ws.Range("M131").Value = -330752.1
$ws.Range("H139").Value = 8041.9546
$ws.Range("I139").Value = 7365.5
$ws.Range("J139").Value = 8428.5
$ws.Range("K139").Value = 22096.5
$ws.Range("L139").Value = 25285.5
$ws.Range("M139").Value = -16956.5
$ws.Range("N139").Value = -35565.5
$ws.Range("H140").Value = 870.375
$ws.Range("I140").Value = 775.7143
$ws.Range("J140").Value = 1533
$ws.Range("K140").Value = 2327.1429
$ws.Range("L140").Value = 4599
$ws.Range("M140").Value = 2852.8571
$ws.Range("N140").Value = -14959
$ws = $wb.Worksheets.Item(6)
$ws.Range("H34").Value = 53333.332
$ws.Range("J34").Value = 53333.332
$ws.Range("L34").Value = 53333.332
$ws.Range("N34").Value = -53869.332
$ws.Range("H35").Value = 20555.5
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 20555.5
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 20555.5
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -21151.5
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 4125502
$ws.Range("I70").Value = 7784504
$ws.Range("J70").Value = 9124.875
$ws.Range("K70").Value = 7784504
$ws.Range("L70").Value = 9124.875
$ws.Range("M70").Value = -7784234
$ws.Range("N70").Value = -9664.875
$ws.Range("H73").Value = 4125502
$ws.Range("I73").Value = 7784504
$ws.Range("J73").Value = 9124.875
$ws.Range("K73").Value = 7784504
$ws.Range("L73").Value = 9124.875
$ws.Range("M73").Value = -7783568
$ws.Range("N73").Value = -10996.875
$ws.Range("H76").Value = 53333.332
$ws.Range("J76").Value = 53333.332
$ws.Range("L76").Value = 53333.332
$ws.Range("N76").Value = -53963.332
$ws.Range("H79").Value = 53333.332
$ws.Range("J79").Value = 53333.332
$ws.Range("L79").Value = 53333.332
$ws.Range("N79").Value = -55517.332
$ws.Range("H97").Value = 523.8889
$ws.Range("I97").Value = 556.2857
$ws.Range("K97").Value = 556.2857
$ws.Range("M97").Value = -60.28570000000002
$ws.Range("H102").Value = 3213.2856
$ws.Range("J102").Value = 2433
$ws.Range("L102").Value = 2433
$ws.Range("N102").Value = -5677
$ws.Range("H113").Value = 1830
$ws.Range("I113").Value = 1830
$ws.Range("K113").Value = 1830
$ws.Range("M113").Value = 340
$ws.Range("H122").Value = 2498.2
$ws.Range("I122").Value = 2498.2
$ws.Range("K122").Value = 7494.599999999999
$ws.Range("M122").Value = -5044.599999999999
$ws.Range("H126").Value = 4075.2727
$ws.Range("I126").Value = 2994
$ws.Range("J126").Value = 4480.75
$ws.Range("K126").Value = 8982
$ws.Range("L126").Value = 13442.25
$ws.Range("M126").Value = -6512
$ws.Range("N126").Value = -18382.25
$ws.Range("H132").Value = 2846.1904
$ws.Range("I132").Value = 2339.7646
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 7019.293799999999
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -4489.293799999999
$ws.Range("N132").Value = -20055.5
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 1347
$ws.Range("I22").Value = 1112.25
$ws.Range("J22").Value = 1581.75
$ws.Range("K22").Value = 1112.25
$ws.Range("L22").Value = 1581.75
$ws.Range("M22").Value = -817.25
$ws.Range("N22").Value = -2171.75
$ws.Range("H27").Value = 1347
$ws.Range("I27").Value = 1112.25
$ws.Range("J27").Value = 1581.75
$ws.Range("K27").Value = 1112.25
$ws.Range("L27").Value = 1581.75
$ws.Range("M27").Value = -1005.25
$ws.Range("N27").Value = -1795.75
$ws.Range("H40").Value = 5115.5713
$ws.Range("I40").Value = 5204.4
$ws.Range("K40").Value = 5204.4
$ws.Range("M40").Value = -5068.4
$ws.Range("H46").Value = 1533.8334
$ws.Range("I46").Value = 1540
$ws.Range("J46").Value = 1523.5555
$ws.Range("K46").Value = 1540
$ws.Range("L46").Value = 1523.5555
$ws.Range("M46").Value = -1352
$ws.Range("N46").Value = -1899.5555
$ws.Range("H61").Value = 2551.5173
$ws.Range("I61").Value = 2230.2856
$ws.Range("J61").Value = 2851.3333
$ws.Range("K61").Value = 2230.2856
$ws.Range("L61").Value = 2851.3333
$ws.Range("M61").Value = -2028.2856
$ws.Range("N61").Value = -3255.3333
$ws.Range("H100").Value = 3135.2856
$ws.Range("I100").Value = 3399.4
$ws.Range("J100").Value = 2475
$ws.Range("K100").Value = 3399.4
$ws.Range("L100").Value = 2475
$ws.Range("M100").Value = -2858.4
$ws.Range("N100").Value = -3557
$ws.Range("H113").Value = 2551.5173
$ws.Range("I113").Value = 2230.2856
$ws.Range("J113").Value = 2851.3333
$ws.Range("K113").Value = 2230.2856
$ws.Range("L113").Value = 2851.3333
$ws.Range("M113").Value = -60.28560000000016
$ws.Range("N113").Value = -7191.3333
$ws.Range("H130").Value = 67500
$ws.Range("J130").Value = 67500
$ws.Range("L130").Value = 67500
$ws.Range("N130").Value = -77540
$ws.Range("H132").Value = 3937.5
$ws.Range("I132").Value = 3937.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11812.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9282.5
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3181.25
$ws.Range("I136").Value = 3054.4614
$ws.Range("K136").Value = 9163.3842
$ws.Range("M136").Value = -6613.3842
$ws = $wb.Worksheets.Item(8)
$ws.Range("H14").Value = 12499.5
$ws.Range("I14").Value = 12499.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 12499.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -12331.5
$ws.Range("N14").ClearContents()
$ws.Range("H41").Value = 14999
$ws.Range("J41").Value = 14999
$ws.Range("L41").Value = 14999
$ws.Range("N41").Value = -15779
$ws.Range("H62").Value = 4140.4287
$ws.Range("I62").Value = 4164
$ws.Range("K62").Value = 4164
$ws.Range("M62").Value = -3540
$ws.Range("H65").Value = 4140.4287
$ws.Range("I65").Value = 4164
$ws.Range("K65").Value = 20820
$ws.Range("M65").Value = -17700
$ws.Range("H81").Value = 1083
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1083
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608
$ws.Range("H113").Value = 999.125
$ws.Range("I113").Value = 639
$ws.Range("K113").Value = 1917
$ws.Range("M113").Value = 253
$ws.Range("H121").Value = 69420
$ws.Range("J121").Value = 69420
$ws.Range("L121").Value = 69420
$ws.Range("N121").Value = -72914
